$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 6).Value = 2.34
$ws.Cells.Item(2, 7).Value = 2.54
$ws.Cells.Item(2, 8).Value = 3.5
$ws.Cells.Item(2, 9).Value = 4.1
$ws.Cells.Item(2, 10).Value = 2.78
$ws.Cells.Item(2, 11).Value = 3.15
$ws.Cells.Item(2, 12).Value = 1.59
$ws.Cells.Item(2, 13).Value = 1.12
$ws.Cells.Item(2, 15).Value = 1.57
$ws.Cells.Item(2, 17).Value = 2.7
$ws.Cells.Item(2, 19).Value = 5.7
$ws.Cells.Item(2, 20).Value = 2.14
$ws.Cells.Item(2, 21).Value = 1.71
$ws.Cells.Item(2, 22).Value = 1.33
$ws.Cells.Item(2, 23).Value = 1.64
$ws.Cells.Item(2, 24).Value = 9.199999999999999
$ws.Cells.Item(2, 25).Value = 11.5
$ws.Cells.Item(2, 26).Value = 30
$ws.Cells.Item(2, 28).Value = 9.199999999999999
$ws.Cells.Item(2, 29).Value = 8.199999999999999
$ws.Cells.Item(2, 30).Value = 19.5
$ws.Cells.Item(2, 32).Value = 17
$ws.Cells.Item(2, 33).Value = 16.5
$ws.Cells.Item(2, 36).Value = 44
$ws.Cells.Item(2, 37).Value = 46
$ws.Cells.Item(2, 40).Value = 48

# Row 3
$ws.Cells.Item(3, 8).Value = 2.06
$ws.Cells.Item(3, 14).Value = 4.6
$ws.Cells.Item(3, 17).Value = 1.64
$ws.Cells.Item(3, 18).Value = 1.52
$ws.Cells.Item(3, 21).Value = 2.36
$ws.Cells.Item(3, 24).Value = 24
$ws.Cells.Item(3, 26).Value = 18.5
$ws.Cells.Item(3, 27).Value = 30
$ws.Cells.Item(3, 28).Value = 19
$ws.Cells.Item(3, 29).Value = 10.5
$ws.Cells.Item(3, 30).Value = 13
$ws.Cells.Item(3, 31).Value = 22
$ws.Cells.Item(3, 33).Value = 16.5
$ws.Cells.Item(3, 34).Value = 18
$ws.Cells.Item(3, 35).Value = 32
$ws.Cells.Item(3, 37).Value = 38
$ws.Cells.Item(3, 39).Value = 75

# Row 4
$ws.Cells.Item(4, 7).Value = 5.7
$ws.Cells.Item(4, 8).Value = 1.82
$ws.Cells.Item(4, 9).Value = 1.95
$ws.Cells.Item(4, 13).Value = 1.09
$ws.Cells.Item(4, 17).Value = 2.16
$ws.Cells.Item(4, 22).Value = 2.02
$ws.Cells.Item(4, 25).Value = 7.6
$ws.Cells.Item(4, 27).Value = 25
$ws.Cells.Item(4, 29).Value = 8.4
$ws.Cells.Item(4, 31).Value = 27
$ws.Cells.Item(4, 37).Value = 1000
$ws.Cells.Item(4, 40).Value = 1000
$ws.Cells.Item(4, 41).Value = 20

# Row 5
$ws.Cells.Item(5, 6).Value = 3.05
$ws.Cells.Item(5, 7).Value = 3.4
$ws.Cells.Item(5, 8).Value = 1.04
$ws.Cells.Item(5, 10).Value = 4.2
$ws.Cells.Item(5, 23).Value = 1.41

# Row 8
$ws.Cells.Item(8, 14).Value = 3.55
$ws.Cells.Item(8, 16).Value = 1.87
$ws.Cells.Item(8, 17).Value = 2.08
$ws.Cells.Item(8, 18).Value = 1.32
$ws.Cells.Item(8, 19).Value = 3.9
$ws.Cells.Item(8, 20).Value = 2.58
$ws.Cells.Item(8, 21).Value = 1.59
$ws.Cells.Item(8, 25).Value = 27
$ws.Cells.Item(8, 32).Value = 6.8
$ws.Cells.Item(8, 34).Value = 42

# Row 9
$ws.Cells.Item(9, 6).Value = 3.05
$ws.Cells.Item(9, 7).Value = 3.1
$ws.Cells.Item(9, 8).Value = 2.6
$ws.Cells.Item(9, 9).Value = 2.62
$ws.Cells.Item(9, 10).Value = 3.4
$ws.Cells.Item(9, 11).Value = 3.45
$ws.Cells.Item(9, 14).Value = 3.75
$ws.Cells.Item(9, 15).Value = 1.34
$ws.Cells.Item(9, 16).Value = 1.91
$ws.Cells.Item(9, 22).Value = 1.61
$ws.Cells.Item(9, 23).Value = 1.47
$ws.Cells.Item(9, 26).Value = 16
$ws.Cells.Item(9, 27).Value = 38
$ws.Cells.Item(9, 28).Value = 12
$ws.Cells.Item(9, 32).Value = 19.5
$ws.Cells.Item(9, 35).Value = 42
$ws.Cells.Item(9, 36).Value = 50
$ws.Cells.Item(9, 41).Value = 25

# Row 10
$ws.Cells.Item(10, 6).Value = 3.6
$ws.Cells.Item(10, 7).Value = 3.65
$ws.Cells.Item(10, 16).Value = 2.12
$ws.Cells.Item(10, 17).Value = 1.87
$ws.Cells.Item(10, 19).Value = 3.15
$ws.Cells.Item(10, 20).Value = 1.72
$ws.Cells.Item(10, 21).Value = 2.32
$ws.Cells.Item(10, 23).Value = 1.37

# Row 11
$ws.Cells.Item(11, 6).Value = 1.88
$ws.Cells.Item(11, 7).Value = 1.89
$ws.Cells.Item(11, 8).Value = 4.9
$ws.Cells.Item(11, 9).Value = 5
$ws.Cells.Item(11, 15).Value = 1.37
$ws.Cells.Item(11, 16).Value = 1.87
$ws.Cells.Item(11, 17).Value = 2.12
$ws.Cells.Item(11, 20).Value = 1.96
$ws.Cells.Item(11, 22).Value = 1.25
$ws.Cells.Item(11, 23).Value = 2.12
$ws.Cells.Item(11, 25).Value = 16
$ws.Cells.Item(11, 27).Value = 120
$ws.Cells.Item(11, 28).Value = 8
$ws.Cells.Item(11, 30).Value = 19
$ws.Cells.Item(11, 31).Value = 70
$ws.Cells.Item(11, 40).Value = 14
$ws.Cells.Item(11, 41).Value = 80

# Row 12
$ws.Cells.Item(12, 6).Value = 3.25
$ws.Cells.Item(12, 7).Value = 3.35
$ws.Cells.Item(12, 8).Value = 2.42
$ws.Cells.Item(12, 9).Value = 2.44
$ws.Cells.Item(12, 10).Value = 3.5
$ws.Cells.Item(12, 11).Value = 3.55
$ws.Cells.Item(12, 13).Value = 1.08
$ws.Cells.Item(12, 18).Value = 1.37
$ws.Cells.Item(12, 19).Value = 3.5
$ws.Cells.Item(12, 20).Value = 1.75
$ws.Cells.Item(12, 22).Value = 1.69
$ws.Cells.Item(12, 23).Value = 1.43
$ws.Cells.Item(12, 25).Value = 11
$ws.Cells.Item(12, 26).Value = 15.5
$ws.Cells.Item(12, 27).Value = 32
$ws.Cells.Item(12, 28).Value = 13.5
$ws.Cells.Item(12, 31).Value = 24
$ws.Cells.Item(12, 32).Value = 23
$ws.Cells.Item(12, 34).Value = 16.5
$ws.Cells.Item(12, 36).Value = 55
$ws.Cells.Item(12, 39).Value = 85
$ws.Cells.Item(12, 41).Value = 19

# Row 13
$ws.Cells.Item(13, 12).Value = 1.37
$ws.Cells.Item(13, 18).Value = 1.44
$ws.Cells.Item(13, 23).Value = 1.79
$ws.Cells.Item(13, 24).Value = 16.5
$ws.Cells.Item(13, 26).Value = 24
$ws.Cells.Item(13, 35).Value = 38
$ws.Cells.Item(13, 38).Value = 30
$ws.Cells.Item(13, 40).Value = 15

# Row 14
$ws.Cells.Item(14, 8).Value = 7.8
$ws.Cells.Item(14, 9).Value = 8
$ws.Cells.Item(14, 16).Value = 3.15
$ws.Cells.Item(14, 18).Value = 1.87
$ws.Cells.Item(14, 19).Value = 2.1
$ws.Cells.Item(14, 21).Value = 2.46
$ws.Cells.Item(14, 24).Value = 36
$ws.Cells.Item(14, 27).Value = 210
$ws.Cells.Item(14, 28).Value = 14
$ws.Cells.Item(14, 36).Value = 14.5
$ws.Cells.Item(14, 40).Value = 4.5
$ws.Cells.Item(14, 41).Value = 65

# Row 15
$ws.Cells.Item(15, 6).Value = 1.85
$ws.Cells.Item(15, 7).Value = 1.86
$ws.Cells.Item(15, 8).Value = 5.1
$ws.Cells.Item(15, 9).Value = 5.3
$ws.Cells.Item(15, 12).Value = 1.47
$ws.Cells.Item(15, 16).Value = 1.83
$ws.Cells.Item(15, 17).Value = 2.16
$ws.Cells.Item(15, 18).Value = 1.31
$ws.Cells.Item(15, 21).Value = 1.94
$ws.Cells.Item(15, 27).Value = 130
$ws.Cells.Item(15, 28).Value = 7.8
$ws.Cells.Item(15, 30).Value = 20
$ws.Cells.Item(15, 32).Value = 10
$ws.Cells.Item(15, 36).Value = 19

# Row 16
$ws.Cells.Item(16, 6).Value = 2.58
$ws.Cells.Item(16, 7).Value = 2.62
$ws.Cells.Item(16, 8).Value = 3.3
$ws.Cells.Item(16, 10).Value = 3.15
$ws.Cells.Item(16, 11).Value = 3.2
$ws.Cells.Item(16, 16).Value = 1.64
$ws.Cells.Item(16, 17).Value = 2.48
$ws.Cells.Item(16, 19).Value = 4.9
$ws.Cells.Item(16, 21).Value = 1.92
$ws.Cells.Item(16, 23).Value = 1.61
$ws.Cells.Item(16, 25).Value = 10
$ws.Cells.Item(16, 28).Value = 8.6
$ws.Cells.Item(16, 32).Value = 14.5
$ws.Cells.Item(16, 35).Value = 70
$ws.Cells.Item(16, 36).Value = 38
$ws.Cells.Item(16, 37).Value = 34
$ws.Cells.Item(16, 38).Value = 60
$ws.Cells.Item(16, 40).Value = 34

# Row 17
$ws.Cells.Item(17, 6).Value = 2.4
$ws.Cells.Item(17, 10).Value = 3.1
$ws.Cells.Item(17, 12).Value = 1.52
$ws.Cells.Item(17, 20).Value = 2.04
$ws.Cells.Item(17, 23).Value = 1.69
$ws.Cells.Item(17, 24).Value = 9.6

# Row 18
$ws.Cells.Item(18, 10).Value = 5.3
$ws.Cells.Item(18, 11).Value = 5.4
$ws.Cells.Item(18, 16).Value = 2.32
$ws.Cells.Item(18, 17).Value = 1.7
$ws.Cells.Item(18, 18).Value = 1.53
$ws.Cells.Item(18, 19).Value = 2.8
$ws.Cells.Item(18, 20).Value = 1.97
$ws.Cells.Item(18, 21).Value = 1.98
$ws.Cells.Item(18, 22).Value = 3.4
$ws.Cells.Item(18, 25).Value = 9
$ws.Cells.Item(18, 26).Value = 8.6
$ws.Cells.Item(18, 27).Value = 12
$ws.Cells.Item(18, 29).Value = 11.5
$ws.Cells.Item(18, 31).Value = 14.5
$ws.Cells.Item(18, 36).Value = 380
$ws.Cells.Item(18, 37).Value = 150
$ws.Cells.Item(18, 38).Value = 130
$ws.Cells.Item(18, 39).Value = 150
$ws.Cells.Item(18, 40).Value = 170
$ws.Cells.Item(18, 41).Value = 5.9

# Row 19
$ws.Cells.Item(19, 8).Value = 1.72
$ws.Cells.Item(19, 9).Value = 1.74
$ws.Cells.Item(19, 10).Value = 4.2
$ws.Cells.Item(19, 11).Value = 4.3
$ws.Cells.Item(19, 20).Value = 1.8
$ws.Cells.Item(19, 22).Value = 2.34
$ws.Cells.Item(19, 40).Value = 70
$ws.Cells.Item(19, 41).Value = 9.199999999999999

# Row 20
$ws.Cells.Item(20, 6).Value = 1.75
$ws.Cells.Item(20, 7).Value = 1.77
$ws.Cells.Item(20, 9).Value = 5.6
$ws.Cells.Item(20, 20).Value = 1.86
$ws.Cells.Item(20, 26).Value = 42
$ws.Cells.Item(20, 30).Value = 20
$ws.Cells.Item(20, 31).Value = 70
$ws.Cells.Item(20, 36).Value = 17.5
$ws.Cells.Item(20, 41).Value = 70
